$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header cell C1: "Item_number" text
$ws.Range("C1").Value = "Item_number"

# Row 2
$ws.Range("B2").Value = 1202136
$ws.Range("C2").Value = 1774013

# Row 3
$ws.Range("B3").Value = 1202136
$ws.Range("C3").Value = 1774014

# Row 4
$ws.Range("B4").Value = 1202136
$ws.Range("C4").Value = 1774015

# Row 5
$ws.Range("B5").Value = 1202136
$ws.Range("C5").Value = 1774016

# Update selection to C5
$ws.Range("C5").Select()
